$d = $word.ActiveDocument

# This "build site" commit removes the trailing footer block that used to
# follow the bibliography entry "Thomson Pioneira (2008)." -- namely the
# blank paragraph, the "Ver no Jupiter ..." paragraph and the
# "(c) 2020 ... Creative Commons Attribution" paragraph. The blank
# paragraph and the page-break paragraph that come after them are left
# untouched.

# Locate the start of the "Ver no Jupiter..." paragraph.
$startFind = $d.Content.Duplicate
$startFind.Find.Execute("Ver no Jupiter Salvar em pdf Salvar em docx", `
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
# Step one character back to also grab the blank paragraph mark that
# precedes it (the paragraph right after "Thomson Pioneira (2008).").
$startPos = $startFind.Start - 1

# Locate the end of the copyright/footer paragraph, including its own
# trailing paragraph mark, so the whole paragraph disappears along with
# everything found above.
$copyrightText = [char]0xA9 + " 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github pages. Original theme under Creative Commons Attribution"
$endFind = $d.Content.Duplicate
$endFind.Find.Execute($copyrightText, `
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$endPos = $endFind.End + 1

$victim = $d.Range($startPos, $endPos)
$victim.Delete()
